# Generate Report for Handoff
# Update the "Latest Handoff" timestamps for the files that were just
# (re-)handed off, across the Overview roll-up sheet and the per-locale
# detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 4, 6, 7, 8, 9, 10 correspond to the files whose handoff was just
# (re-)generated; row 2/3 (already handed back) and row 5 (still in
# translation) are untouched.
$rows = @(4, 6, 7, 8, 9, 10)

foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-03-19 17:29:54"
    $zhcn.Range("E$r").Value     = "2016-03-19 17:29:45"
    $dede.Range("E$r").Value     = "2016-03-19 17:29:54"
}
